$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# Sheet "position" (1st sheet): zoom 100 -> 150
# ---------------------------------------------------------------------------
$ws1 = $wb.Worksheets.Item(1)
$ws1.Activate()
$excel.ActiveWindow.Zoom = 150

# ---------------------------------------------------------------------------
# Sheet "e2Single" (2nd sheet): no longer the tab-selected sheet,
# zoom 91 -> 150, selection moves from K13 to N24, view scrolls to A9
# ---------------------------------------------------------------------------
$ws2 = $wb.Worksheets.Item(2)
$ws2.Activate()
$excel.ActiveWindow.Zoom = 150
$excel.ActiveWindow.ScrollRow = 9
$excel.ActiveWindow.ScrollColumn = 1
$ws2.Range("N24").Select()

# ---------------------------------------------------------------------------
# New sheet "listBindTest" appended after "e2Single"
# ---------------------------------------------------------------------------
$ws3 = $wb.Worksheets.Add($null, $ws2)
$ws3.Name = "listBindTest"

# Header / label cells (plain strings -> new shared-string entries T1, T2, t3)
$ws3.Range("A1").Value = "T1"
$ws3.Range("A2").Value = "T2"
$ws3.Range("A3").Value = "t3"

# C4
$ws3.Range("C4").Value = 1

# Row 5 - all ones
$ws3.Range("A5").Value = 1
$ws3.Range("B5").Value = 1
$ws3.Range("C5").Value = 1
$ws3.Range("D5").Value = 1
$ws3.Range("E5").Value = 1
$ws3.Range("F5").Value = 1

# Row 6
$ws3.Range("A6").Value = 1
$ws3.Range("B6").Value = 1
$ws3.Range("C6").Value = 1
$ws3.Range("D6").Value = 11
$ws3.Range("E6").Value = 16
$ws3.Range("F6").Value = 1

# Row 7
$ws3.Range("A7").Value = 1
$ws3.Range("B7").Value = 2
$ws3.Range("C7").Value = 7
$ws3.Range("D7").Value = 12
$ws3.Range("E7").Value = 17
$ws3.Range("F7").Value = 1

# Row 8 (D8 stays empty but keeps the "orange" highlight)
$ws3.Range("A8").Value = 1
$ws3.Range("B8").Value = 3
$ws3.Range("C8").Value = 8
$ws3.Range("E8").Value = 18
$ws3.Range("F8").Value = 1

# Row 9 (E9 stays empty but keeps the "yellow" highlight)
$ws3.Range("A9").Value = 1
$ws3.Range("B9").Value = 4
$ws3.Range("C9").Value = 9
$ws3.Range("D9").Value = 14
$ws3.Range("F9").Value = 1
# (E9 intentionally left blank)

# Row 10
$ws3.Range("A10").Value = 1
$ws3.Range("B10").Value = 5
$ws3.Range("C10").Value = 10
$ws3.Range("D10").Value = 15
$ws3.Range("E10").Value = 20
$ws3.Range("F10").Value = 1

# Row 11 (no A11)
$ws3.Range("B11").Value = 1
$ws3.Range("C11").Value = 1
$ws3.Range("D11").Value = 2
$ws3.Range("E11").Value = 1
$ws3.Range("F11").Value = 1

# Row 12 (no A12, no B12)
$ws3.Range("C12").Value = 1
$ws3.Range("D12").Value = 4
$ws3.Range("E12").Value = 1
$ws3.Range("F12").Value = 1

# ---------------------------------------------------------------------------
# Highlight fills: yellow block B6:E10, except D8 which is orange (and left
# empty) and E9 which is yellow but left empty.
# ---------------------------------------------------------------------------
$yellowRange = $ws3.Range("B6:E10")
$yellowRange.Interior.Color = 65535
$yellowRange.Interior.PatternColor = 65535

$orangeRange = $ws3.Range("D8")
$orangeRange.Interior.Color = 16639
$orangeRange.Interior.PatternColor = 255

# ---------------------------------------------------------------------------
# Sheet "listBindTest" view: zoom 150, tab-selected, selection E18
# ---------------------------------------------------------------------------
$ws3.Activate()
$excel.ActiveWindow.Zoom = 150
$ws3.Range("E18").Select()

Write-Host "done"
